$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells F1:H1, matching style of existing header (A1:E1 use style "1")
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

# Copy the header style (bold, centered, bordered) from an existing header cell
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Boolean values for new columns F (KNN_Outliers_MAD), G (SVM_Outliers_MAD), H (RF_Outliers_MAD)
$values = @{
    2  = @($false, $false, $false)
    3  = @($false, $false, $false)
    4  = @($false, $false, $false)
    5  = @($false, $false, $false)
    6  = @($false, $false, $false)
    7  = @($false, $false, $false)
    8  = @($false, $true,  $false)
    9  = @($false, $false, $false)
    10 = @($false, $false, $false)
    11 = @($false, $false, $false)
    12 = @($false, $false, $false)
    13 = @($false, $false, $false)
    14 = @($false, $false, $false)
    15 = @($false, $false, $false)
    16 = @($false, $false, $false)
    17 = @($false, $false, $false)
    18 = @($false, $false, $false)
    19 = @($false, $false, $false)
    20 = @($true,  $true,  $true)
    21 = @($false, $false, $false)
}

foreach ($row in $values.Keys | Sort-Object) {
    $v = $values[$row]
    $ws.Cells.Item($row, 6).Value = $v[0]
    $ws.Cells.Item($row, 7).Value = $v[1]
    $ws.Cells.Item($row, 8).Value = $v[2]
}

$wb.Save()
